$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2 through 303 contain a date serial value (45202) that
# needs to be bumped to 45203 (i.e. one day later), keeping existing
# number formatting/style intact.
$ws.Range("C2:C303").Value2 = 45203
